$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.955.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.353.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  +3.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.86"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.22%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.53%  "
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.108"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.705.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.909"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.330.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.868.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  -2.32%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  +15.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.137"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.33%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.127"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.53%  "
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0753"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("E36").Value = "  +4.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  +5.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.214"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.45%  "
$ws.Range("E42").Value = "  +11.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.73%  "
$ws.Range("E47").Value = "  +8.32%  "
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "100.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.71%  "
